# Update "人气"/popularity counts (column F) on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 9818
$wsExhibit.Range("F7").Value = 889
$wsExhibit.Range("F10").Value = 3900
$wsExhibit.Range("F11").Value = 170
$wsExhibit.Range("F13").Value = 40
$wsExhibit.Range("F16").Value = 541
$wsExhibit.Range("F18").Value = 267
$wsExhibit.Range("F19").Value = 1440

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 9818
$wsAll.Range("F8").Value = 889
$wsAll.Range("F11").Value = 3900
$wsAll.Range("F12").Value = 170
$wsAll.Range("F14").Value = 40
$wsAll.Range("F17").Value = 541
$wsAll.Range("F19").Value = 267
$wsAll.Range("F20").Value = 1440
